$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price-like text (e.g. "87.899.20", "0.386") that Excel would
# otherwise auto-coerce to a number when assigned via .Value. Force the column
# to Text format while writing, then restore the default (un-styled) look so the
# cells keep their original "no explicit style" appearance.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '87.899.20'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '3.246.04'
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '207.53'
$ws.Range("E5").Value = '  -4.45%  '
$ws.Range("D6").Value = '618.26'
$ws.Range("E6").Value = '  -4.63%  '
$ws.Range("D7").Value = '0.386'
$ws.Range("E7").Value = '  +10.49%  '
$ws.Range("D8").Value = '0.689'
$ws.Range("E8").Value = '  +14.42%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("D10").Value = '3.244.86'
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("D11").Value = '0.552'
$ws.Range("E11").Value = '  -5.43%  '
$ws.Range("D12").Value = '0.184'
$ws.Range("E12").Value = '  +9.49%  '
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value = '  -5.13%  '
$ws.Range("D14").Value = '5.43'
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").Value = '3.837.95'
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").Value = '33.18'
$ws.Range("E16").Value = '  -5.89%  '
$ws.Range("D17").Value = '87.604.21'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '3.211.69'
$ws.Range("E18").Value = '  -2.61%  '
$ws.Range("D19").Value = '13.72'
$ws.Range("E19").Value = '  -6.41%  '
$ws.Range("E20").Value = '  -4.24%  '
$ws.Range("D21").Value = '427.09'
$ws.Range("E21").Value = '  -6.10%  '
$ws.Range("D22").Value = '8.75'
$ws.Range("E22").Value = '  -14.43%  '
$ws.Range("D23").Value = '5.20'
$ws.Range("E23").Value = '  -6.05%  '
$ws.Range("D24").Value = '5.26'
$ws.Range("E24").Value = '  -4.20%  '
$ws.Range("D25").Value = '11.88'
$ws.Range("E25").Value = '  -5.24%  '
$ws.Range("D26").Value = '3.399.00'
$ws.Range("E26").Value = '  -2.76%  '
$ws.Range("D27").Value = '0.0000137'
$ws.Range("E27").Value = '  +8.82%  '
$ws.Range("D28").Value = '75.18'
$ws.Range("E28").Value = '  -4.18%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").Value = '0.179'
$ws.Range("E30").Value = '  -8.56%  '
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").Value = '8.58'
$ws.Range("E32").Value = '  -8.20%  '
$ws.Range("D33").Value = '551.28'
$ws.Range("E33").Value = '  -9.20%  '
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  -8.52%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '1.30'
$ws.Range("E35").Value = '  -19.16%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '6.82'
$ws.Range("E36").Value = '  -4.30%  '
$ws.Range("D37").Value = '0.137'
$ws.Range("E37").Value = '  -5.84%  '
$ws.Range("D38").Value = '22.59'
$ws.Range("E38").Value = '  -3.57%  '
$ws.Range("D39").Value = '21.86'
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").Value = '3.06'
$ws.Range("E41").Value = '  +2.20%  '
$ws.Range("D42").Value = '0.391'
$ws.Range("E42").Value = '  -6.54%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '1.94'
$ws.Range("E43").Value = '  -10.57%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '150.20'
$ws.Range("E45").Value = '  -5.47%  '
$ws.Range("D46").Value = '176.77'
$ws.Range("E46").Value = '  -7.88%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = '0.131'
$ws.Range("E47").Value = '  +14.89%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").Value = '43.95'
$ws.Range("E48").Value = '  -4.63%  '
$ws.Range("D49").Value = '1.31'
$ws.Range("E49").Value = '  -6.83%  '
$ws.Range("D50").Value = '4.13'
$ws.Range("E50").Value = '  -7.34%  '
$ws.Range("D51").Value = '0.616'
$ws.Range("E51").Value = '  -7.03%  '

$dRange.Style = "Normal"
